# Applies the cryptos-list price/volume refresh described in the commit
# 'Updated cryptos list ... with GitHub Actions'.
# Numeric-looking Price (column D) values are apostrophe-prefixed so
# Excel stores them as text (matching the workbook's existing inlineStr
# convention) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.994.63"
$ws.Range("E2").Value = "  +2.09%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.297.91"
$ws.Range("E3").Value = "  +1.50%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'301.55"
$ws.Range("E5").Value = "  +1.06%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'98.62"
$ws.Range("E6").Value = "  +4.83%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.97%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.506"
$ws.Range("E9").Value = "  +3.33%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'34.13"
$ws.Range("E10").Value = "  +3.32%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0798"
$ws.Range("E11").Value = "  +0.95%  "

# Row 12 - OKB
$ws.Range("D12").Value = "'49.03"
$ws.Range("E12").Value = "  +2.18%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +4.24%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'17.86"
$ws.Range("E14").Value = "  +15.35%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +1.95%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.654.07"
$ws.Range("E16").Value = "  +1.46%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.291.77"
$ws.Range("E17").Value = "  +0.77%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +4.63%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "42.908.22"
$ws.Range("E19").Value = "  +1.91%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.35"
$ws.Range("E20").Value = "  +8.84%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").Value = "  +1.82%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.46%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'67.70"

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "'236.35"
$ws.Range("E24").Value = "  +1.42%  "

# Row 25 - ImmutableX
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  +13.30%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.08%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  +0.68%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'24.72"
$ws.Range("E28").Value = "  +3.95%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -4.90%  "

# Row 30 - Monero
$ws.Range("D30").Value = "'167.41"
$ws.Range("E30").Value = "  -0.10%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "'33.76"
$ws.Range("E31").Value = "  +0.38%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  +1.13%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.00%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.05"
$ws.Range("E34").Value = "  +2.65%  "

# Row 35 - RenderToken
$ws.Range("D35").Value = "'4.55"
$ws.Range("E35").Value = "  +1.66%  "

# Row 36 - WEMIXToken
$ws.Range("D36").Value = "'2.41"
$ws.Range("E36").Value = "  +3.33%  "

# Row 37 - Celestia
$ws.Range("D37").Value = "'16.88"
$ws.Range("E37").Value = "  +4.58%  "

# Row 38 - Hedera
$ws.Range("D38").Value = "'0.0691"
$ws.Range("E38").Value = "  +0.26%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +3.96%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +4.77%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  +0.74%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +0.25%  "

# Row 43 - ApeXProtocol
$ws.Range("E43").Value = "  -2.92%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.989.49"
$ws.Range("E44").Value = "  +1.58%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +2.50%  "

# Row 46 - FraxShare
$ws.Range("D46").Value = "'9.98"
$ws.Range("E46").Value = "  +5.02%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "'17.58"
$ws.Range("E47").Value = "  +1.60%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  +2.78%  "

# Row 49 - MultiversX
$ws.Range("D49").Value = "'56.28"
$ws.Range("E49").Value = "  +9.17%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.525.37"
$ws.Range("E50").Value = "  +1.45%  "

# Row 51 - Stacks
$ws.Range("E51").Value = "  +3.61%  "

